$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the used range down to row 551, copying the formatting of the last
# existing data row (469) so the new rows inherit the same style (s="1" on col A).
$ws.Range("A469:F469").Copy()
$ws.Range("A470:F551").PasteSpecial(-4122)

$dates = @("02/ago/2021","03/ago/2021","04/ago/2021","05/ago/2021","06/ago/2021","07/ago/2021","08/ago/2021","09/ago/2021","10/ago/2021","11/ago/2021","12/ago/2021","13/ago/2021","14/ago/2021","15/ago/2021","16/ago/2021","17/ago/2021","18/ago/2021","19/ago/2021","20/ago/2021","21/ago/2021","22/ago/2021","23/ago/2021","24/ago/2021","25/ago/2021","26/ago/2021","27/ago/2021","28/ago/2021","29/ago/2021","30/ago/2021","31/ago/2021","01/set/2021","02/set/2021","03/set/2021","04/set/2021","05/set/2021","06/set/2021","07/set/2021","08/set/2021","09/set/2021","10/set/2021","11/set/2021","12/set/2021","13/set/2021","14/set/2021","15/set/2021","16/set/2021","17/set/2021","18/set/2021","19/set/2021","20/set/2021","21/set/2021","22/set/2021","23/set/2021","24/set/2021","25/set/2021","26/set/2021","27/set/2021","28/set/2021","29/set/2021","30/set/2021","01/out/2021","02/out/2021","03/out/2021","04/out/2021","05/out/2021","06/out/2021","07/out/2021","08/out/2021","09/out/2021","10/out/2021","11/out/2021","12/out/2021","13/out/2021","14/out/2021","15/out/2021","16/out/2021","17/out/2021","18/out/2021","19/out/2021","20/out/2021","21/out/2021","22/out/2021")
$casos = @(8703,8703,8703,8703,8703,8703,8703,8703,8703,8703,8703,8750,8750,8755,8755,8758,8761,8761,8761,8766,8771,8771,8773,8773,8786,8786,8786,8786,8786,8786,8796,8798,8798,8798,8798,8798,8818,8818,8818,8818,8818,8818,8824,8824,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8825,8839,8840,8845,8845,8845,8850,8851,8851,8851,8854,8854,8854,8854,8862,8868,8872,8881,8881)
$mortos = @(223,223,223,223,223,223,223,223,223,223,223,229,229,229,229,229,230,230,230,232,232,232,232,232,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,233,234,234,234,234,234,235,235,235,235,235,235,235,235,235,235,235,235,235)
$mortesDiarias = @(0,0,0,0,0,0,0,0,0,0,0,6,0,0,0,0,1,0,0,2,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0)
$casosDiarios = @(0,0,0,0,0,0,0,0,0,0,0,47,0,5,0,3,3,0,0,5,5,0,2,0,13,0,0,0,0,0,10,2,0,0,0,0,20,0,0,0,0,0,6,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,14,1,5,0,0,5,1,0,0,3,0,0,0,8,6,4,9,0)
$population = 118516

$startRow = 470
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $casos[$i]
    $ws.Cells.Item($r, 3).Value = $mortos[$i]
    $ws.Cells.Item($r, 4).Value = $mortesDiarias[$i]
    $ws.Cells.Item($r, 5).Value = $casosDiarios[$i]
    $ws.Cells.Item($r, 6).Value = $population
}

Write-Output "Added rows 470 to 551"
